# Links_youtube.xlsx - add a new "Practice tasks 3, 4, 5, and 6" row with its
# YouTube link (row 10), matching the row-6 hyperlink style, then leave the
# selection/page setup the way the authored workbook ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row. Write the URL first, then the label, so new shared-string
# entries land in the same order as the authored file (URL, then label).
$ws.Range("B10").Value = "https://youtu.be/NtVpk2mr_kM"
$ws.Range("A10").Value = "Practice tasks 3, 4, 5, and 6"

# Turn B10 into a real hyperlink, like B6.
$ws.Hyperlinks.Add($ws.Range("B10"), "https://youtu.be/NtVpk2mr_kM")

# Give B10 the same "Hyperlink" cell style used by the other hyperlink
# cell (B6), instead of the default cell style.
$ws.Range("B10").Style = "Гиперссылка"

# Match the final selection left behind in the saved file.
$ws.Range("D14").Select()

# Page setup the sheet picked up (paper size + orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
